# Fix availability zone 2 -> 1
# On each slide, the shape named "Rectangle 75" (the left-hand
# "Availability Zone" label box) is mislabeled "Availability Zone 2"
# and should read "Availability Zone 1". The other "Availability Zone 2"
# box (shape "Rectangle 8") is left untouched.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $sh = $s.Shapes.Item("Rectangle 75")
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "Availability Zone 2") {
            $sh.TextFrame.TextRange.Text = "Availability Zone 1"
        }
    }
}
